# Insert a new weekly data row at row 205 (pushing the existing rows 205..266
# down to 206..267) and populate it with the new "Ají" price-report entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(205).Insert()

$ws.Cells.Item(205, 1).Value  = 9
$ws.Cells.Item(205, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(205, 3).Value  = "Metropolitana"
$ws.Cells.Item(205, 4).Value  = 44642
$ws.Cells.Item(205, 5).Value  = 13
$ws.Cells.Item(205, 6).Value  = 100112021
$ws.Cells.Item(205, 7).Value  = "Ají"
$ws.Cells.Item(205, 8).Value  = "Americana (o)"
$ws.Cells.Item(205, 9).Value  = "Primera"
$ws.Cells.Item(205, 10).Value = 25
$ws.Cells.Item(205, 11).Value = 22000
$ws.Cells.Item(205, 12).Value = 24000
$ws.Cells.Item(205, 13).Value = 22960
$ws.Cells.Item(205, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(205, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(205, 16).Value = 918
$ws.Cells.Item(205, 17).Value = 25
$ws.Cells.Item(205, 18).Value = "Hortaliza"
